$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131, shifting existing rows 131-214 down to 132-215.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new record's data.
$ws.Range("A131").Value = 8
$ws.Range("B131").Value = "Terminal La Palmera de La Serena"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 45176
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 100112052
$ws.Range("G131").Value = "Albahaca"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 800
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 3500
$ws.Range("M131").Value = 3250
$ws.Range("N131").Value = "`$/paquete"
$ws.Range("O131").Value = "Región de Arica y Parinacota"
$ws.Range("P131").Value = 3250
$ws.Range("Q131").Value = 1
$ws.Range("R131").Value = "Hortaliza"

# Match the date style used by the rest of column D.
$ws.Range("D131").NumberFormat = $ws.Range("D132").NumberFormat
